$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "GENERAL / GENERAL" row (FAMILIA=GENERAL, CATEGORIA=GENERAL)
$ws.Rows.Item(71).Delete() | Out-Null

# Delete the six "SUMINISTROS" family rows (shifted up to 145:150 after the
# previous deletion)
$ws.Range("A145:D150").EntireRow.Delete() | Out-Null

# Re-apply the AutoFilter bookkeeping name over the now-smaller table range
# without leaving a visible autoFilter marker on the sheet.
$n = $ws.Names.Add("_xlnm._FilterDatabase", "=Hoja1!`$A`$1:`$D`$154")
$n.Visible = $false

# Move the active selection
$ws.Range("C6").Select() | Out-Null
